$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 500
$ws1.Range("F6").Value = 390
$ws1.Range("F7").Value = 478
$ws1.Range("F8").Value = 890
$ws1.Range("F9").Value = 117
$ws1.Range("F10").Value = 826
$ws1.Range("F11").Value = 650
$ws1.Range("F12").Value = 122
$ws1.Range("F13").Value = 40
$ws1.Range("F14").Value = 54
$ws1.Range("F16").Value = 216
$ws1.Range("F19").Value = 1232
$ws1.Range("F21").Value = 949
$ws1.Range("F22").Value = 2696
$ws1.Range("F23").Value = 1162
$ws1.Range("F24").Value = 615
$ws1.Range("F26").Value = 1202
$ws1.Range("F30").Value = 1230

# Sheet: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 486

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 500
$ws4.Range("F8").Value = 390
$ws4.Range("F9").Value = 478
$ws4.Range("F10").Value = 486
$ws4.Range("F11").Value = 486
$ws4.Range("F15").Value = 890
$ws4.Range("F16").Value = 117
$ws4.Range("F17").Value = 826
$ws4.Range("F18").Value = 650
$ws4.Range("F19").Value = 122
$ws4.Range("F21").Value = 40
$ws4.Range("F25").Value = 54
$ws4.Range("F28").Value = 216
$ws4.Range("F31").Value = 1232
$ws4.Range("F33").Value = 949
$ws4.Range("F34").Value = 2696
$ws4.Range("F35").Value = 1162
$ws4.Range("F36").Value = 615
$ws4.Range("F38").Value = 1202
$ws4.Range("F43").Value = 1230
